$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 17023.334
$ws.Range("I11").Value = 17023.334
$ws.Range("K11").Value = 17023.334
$ws.Range("M11").Value = -16883.334
$ws.Range("H28").Value = 1624.6
$ws.Range("I28").Value = 1875
$ws.Range("J28").Value = 1249
$ws.Range("K28").Value = 1875
$ws.Range("L28").Value = 1249
$ws.Range("M28").Value = -1390
$ws.Range("N28").Value = -2219
$ws.Range("H116").Value = 1763.6364
$ws.Range("I116").Value = 1542.8572
$ws.Range("J116").Value = 2150
$ws.Range("K116").Value = 1542.8572
$ws.Range("L116").Value = 2150
$ws.Range("M116").Value = 1899.1428
$ws.Range("N116").Value = -9034
$ws.Range("H132").Value = 629968.25
$ws.Range("I132").Value = 1588.5968
$ws.Range("K132").Value = 4765.7904
$ws.Range("M132").Value = -2235.7904
$ws.Range("H133").Value = 60780
$ws.Range("J133").Value = 60780
$ws.Range("L133").Value = 60780
$ws.Range("N133").Value = -70900
$ws.Range("H137").Value = 1924617.5
$ws.Range("I137").Value = 2326579.2
$ws.Range("K137").Value = 6979737.600000001
$ws.Range("M137").Value = -6977187.600000001
$ws.Range("H138").Value = 2284772
$ws.Range("I138").Value = 1230.1111
$ws.Range("K138").Value = 3690.3333
$ws.Range("M138").Value = 1449.6667
$ws.Range("H139").Value = 48000
$ws.Range("J139").Value = 48000
$ws.Range("L139").Value = 48000
$ws.Range("N139").Value = -58280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1447.0571
$ws.Range("I2").Value = 993.6539
$ws.Range("J2").Value = 2756.889
$ws.Range("K2").Value = 993.6539
$ws.Range("L2").Value = 2756.889
$ws.Range("M2").Value = -880.6539
$ws.Range("N2").Value = -2982.889
$ws.Range("H32").Value = 20351.791
$ws.Range("I32").Value = 18776.305
$ws.Range("J32").Value = 25293.092
$ws.Range("K32").Value = 18776.305
$ws.Range("L32").Value = 25293.092
$ws.Range("M32").Value = -18489.305
$ws.Range("N32").Value = -25867.092
$ws.Range("H61").Value = 17580102
$ws.Range("I61").Value = 19628382
$ws.Range("J61").Value = 169715.67
$ws.Range("K61").Value = 19628382
$ws.Range("L61").Value = 169715.67
$ws.Range("M61").Value = -19628170
$ws.Range("N61").Value = -170139.67
$ws.Range("H74").Value = 8180152.5
$ws.Range("I74").Value = 10132371
$ws.Range("J74").Value = 127251.375
$ws.Range("K74").Value = 10132371
$ws.Range("L74").Value = 127251.375
$ws.Range("M74").Value = -10131497
$ws.Range("N74").Value = -128999.375
$ws.Range("H77").Value = 8180152.5
$ws.Range("I77").Value = 10132371
$ws.Range("J77").Value = 127251.375
$ws.Range("K77").Value = 50661855
$ws.Range("L77").Value = 636256.875
$ws.Range("M77").Value = -50657487
$ws.Range("N77").Value = -644992.875
$ws.Range("H110").Value = 303913.25
$ws.Range("I110").Value = 435388.25
$ws.Range("K110").Value = 435388.25
$ws.Range("M110").Value = -433343.25
$ws.Range("H116").Value = 1447.0571
$ws.Range("I116").Value = 993.6539
$ws.Range("J116").Value = 2756.889
$ws.Range("K116").Value = 993.6539
$ws.Range("L116").Value = 2756.889
$ws.Range("M116").Value = 1300.3461
$ws.Range("N116").Value = -7344.889
$ws.Range("H136").Value = 17580102
$ws.Range("I136").Value = 19628382
$ws.Range("J136").Value = 169715.67
$ws.Range("K136").Value = 58885146
$ws.Range("L136").Value = 509147.01
$ws.Range("M136").Value = -58882596
$ws.Range("N136").Value = -514247.01

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1447.0571
$ws.Range("I3").Value = 993.6539
$ws.Range("J3").Value = 2756.889
$ws.Range("K3").Value = 993.6539
$ws.Range("L3").Value = 2756.889
$ws.Range("M3").Value = -879.6539
$ws.Range("N3").Value = -2984.889
$ws.Range("H134").Value = 2276.639
$ws.Range("I134").Value = 1268.0416
$ws.Range("J134").Value = 4293.8335
$ws.Range("K134").Value = 3804.1248
$ws.Range("L134").Value = 12881.5005
$ws.Range("M134").Value = -1269.1248
$ws.Range("N134").Value = -17951.5005
$ws.Range("H138").Value = 32044.445
$ws.Range("J138").Value = 32044.445
$ws.Range("L138").Value = 32044.445
$ws.Range("N138").Value = -42324.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2708.775
$ws.Range("I31").Value = 1301.8
$ws.Range("J31").Value = 5053.7334
$ws.Range("K31").Value = 1301.8
$ws.Range("L31").Value = 5053.7334
$ws.Range("M31").Value = -1006.8
$ws.Range("N31").Value = -5643.7334
$ws.Range("H34").Value = 2708.775
$ws.Range("I34").Value = 1301.8
$ws.Range("J34").Value = 5053.7334
$ws.Range("K34").Value = 1301.8
$ws.Range("L34").Value = 5053.7334
$ws.Range("M34").Value = -1099.8
$ws.Range("N34").Value = -5457.7334
$ws.Range("H58").Value = 18183380
$ws.Range("I58").Value = 23257288
$ws.Range("J58").Value = 1873.5834
$ws.Range("K58").Value = 23257288
$ws.Range("L58").Value = 1873.5834
$ws.Range("M58").Value = -23257085
$ws.Range("N58").Value = -2279.5834
$ws.Range("H132").Value = 26505.9
$ws.Range("I132").Value = 1172.7858
$ws.Range("J132").Value = 85616.5
$ws.Range("K132").Value = 3518.3574
$ws.Range("L132").Value = 256849.5
$ws.Range("M132").Value = -988.3574000000003
$ws.Range("N132").Value = -261909.5
$ws.Range("H134").Value = 27033
$ws.Range("I134").Value = 1317.3636
$ws.Range("J134").Value = 121323.664
$ws.Range("K134").Value = 3952.0908
$ws.Range("L134").Value = 363970.992
$ws.Range("M134").Value = -1417.0908
$ws.Range("N134").Value = -369040.992
$ws.Range("H136").Value = 18183380
$ws.Range("I136").Value = 23257288
$ws.Range("J136").Value = 1873.5834
$ws.Range("K136").Value = 69771864
$ws.Range("L136").Value = 5620.7502
$ws.Range("M136").Value = -69769314
$ws.Range("N136").Value = -10720.7502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 2614.4285
$ws.Range("J76").Value = 3580
$ws.Range("L76").Value = 10740
$ws.Range("N76").Value = -11506
$ws.Range("H79").Value = 2614.4285
$ws.Range("J79").Value = 3580
$ws.Range("L79").Value = 10740
$ws.Range("N79").Value = -13392
$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H94").Value = 1944.8
$ws.Range("I94").Value = 908
$ws.Range("J94").Value = 3500
$ws.Range("K94").Value = 2724
$ws.Range("L94").Value = 10500
$ws.Range("M94").Value = -2048
$ws.Range("N94").Value = -11852
$ws.Range("H100").Value = 2840
$ws.Range("J100").Value = 2840
$ws.Range("L100").Value = 8520
$ws.Range("N100").Value = -10142
$ws.Range("H106").Value = 2700.0417
$ws.Range("J106").Value = 2700.0417
$ws.Range("L106").Value = 8100.125100000001
$ws.Range("N106").Value = -9992.125100000001
$ws.Range("H109").Value = 3509.0908
$ws.Range("I109").Value = 1100
$ws.Range("J109").Value = 3750
$ws.Range("K109").Value = 3300
$ws.Range("L109").Value = 11250
$ws.Range("M109").Value = -2260
$ws.Range("N109").Value = -13330
$ws.Range("H112").Value = 19611332
$ws.Range("I112").Value = 2317.8333
$ws.Range("J112").Value = 30307158
$ws.Range("K112").Value = 6953.499899999999
$ws.Range("L112").Value = 90921474
$ws.Range("M112").Value = -5845.499899999999
$ws.Range("N112").Value = -90923690
$ws.Range("H131").Value = 1471.6097
$ws.Range("J131").Value = 1623.2222
$ws.Range("L131").Value = 4869.6666
$ws.Range("N131").Value = -14949.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1753.75
$ws.Range("I122").Value = 1553.5
$ws.Range("J122").Value = 1954
$ws.Range("K122").Value = 4660.5
$ws.Range("L122").Value = 5862
$ws.Range("M122").Value = -2210.5
$ws.Range("N122").Value = -10762
$ws.Range("H135").Value = 37512.25
$ws.Range("I135").Value = 30709
$ws.Range("J135").Value = 39780
$ws.Range("K135").Value = 30709
$ws.Range("L135").Value = 39780
$ws.Range("M135").Value = -25639
$ws.Range("N135").Value = -49920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4313
$ws.Range("I7").Value = 3048
$ws.Range("J7").Value = 5999.6665
$ws.Range("K7").Value = 3048
$ws.Range("L7").Value = 5999.6665
$ws.Range("M7").Value = -2936
$ws.Range("N7").Value = -6223.6665
$ws.Range("H40").Value = 2870.0286
$ws.Range("I40").Value = 2337.5356
$ws.Range("K40").Value = 2337.5356
$ws.Range("M40").Value = -2201.5356
$ws.Range("H126").Value = 4313
$ws.Range("I126").Value = 3048
$ws.Range("J126").Value = 5999.6665
$ws.Range("K126").Value = 9144
$ws.Range("L126").Value = 17998.9995
$ws.Range("M126").Value = -6674
$ws.Range("N126").Value = -22938.9995
$ws.Range("H132").Value = 36809.277
$ws.Range("I132").Value = 1878.2632
$ws.Range("K132").Value = 5634.7896
$ws.Range("M132").Value = -3104.7896
$ws.Range("H136").Value = 47774.723
$ws.Range("I136").Value = 32256.969
$ws.Range("J136").Value = 92917.27
$ws.Range("K136").Value = 96770.90700000001
$ws.Range("L136").Value = 278751.81
$ws.Range("M136").Value = -94220.90700000001
$ws.Range("N136").Value = -283851.81

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 34867.934
$ws.Range("I132").Value = 21955.213
$ws.Range("J132").Value = 85442.75
$ws.Range("K132").Value = 65865.639
$ws.Range("L132").Value = 256328.25
$ws.Range("M132").Value = -63335.639
$ws.Range("N132").Value = -261388.25
$ws.Range("H136").Value = 44521.617
$ws.Range("I136").Value = 30979.121
$ws.Range("J136").Value = 76443.21000000001
$ws.Range("K136").Value = 92937.363
$ws.Range("L136").Value = 229329.63
$ws.Range("M136").Value = -90387.363
$ws.Range("N136").Value = -234429.63
